# EST-1943: format the stuff that needs to be updated manually in red.
#
# Three placeholder snippets in the Italian order-confirmation template get
# split out into their own run and coloured red (FF0000) so the person who
# fills in the template can find them quickly:
#   1. "DD MMM YYYY"              (offer date, in the opening paragraph)
#   2. "DD MMM YYYY"              (works deadline date)
#   3. "€ X.XXX,00 (XXXX/00)"     (fee amount)

$d = $word.ActiveDocument
$wdColorRed = 255   # OLE/VBA colour value -> 0x0000FF (BGR) -> RGB(255,0,0) -> w:color FF0000

# ---------------------------------------------------------------------
# 1) "Con la presente ... offerta del DD MMM YYYY, Vi confermiamo l'ordine..."
#    The paragraph contains Word's auto "_GoBack" bookmark sitting right
#    before "ordine"; relocate it to sit right after "...nuova of" (i.e.
#    right where the run boundary needs to be) before colouring the date,
#    mirroring how the author's edit reshaped the runs around it.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Con la presente e in riferimento alla Vostra nuova of", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $rng.End
$d.Range($splitPoint, $splitPoint).Bookmarks.Add("_GoBack")

$dateRng = $d.Content
$dateRng.Find.Execute("DD MMM YYYY, Vi confermiamo", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRng.Start
$dateEnd = $dateStart + "DD MMM YYYY".Length
$d.Range($dateStart, $dateEnd).Font.Color = $wdColorRed

# Re-write the trailing text (identical content) so the two runs that now
# sit side by side after the (relocated) bookmark - both with the same
# formatting - collapse back into a single run, exactly as Word would do
# when the text is retyped.
$tailRng = $d.Content
$tailRng.Find.Execute("l’ordine come di seguito precisato.", $false, $false, $false, $false, $false, $true, 1, $false, "l’ordine come di seguito precisato.", 2)

# ---------------------------------------------------------------------
# 2) "I lavori dovranno essere effettuati entro il DD MMM YYYY."
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("I lavori dovranno essere effettuati entro il DD MMM YYYY", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2End = $rng2.End
$rng2Start = $rng2End - "DD MMM YYYY".Length
$d.Range($rng2Start, $rng2End).Font.Color = $wdColorRed

# ---------------------------------------------------------------------
# 3) "L'importo dell'incarico a Voi affidato ammonta a € X.XXX,00 (XXXX/00) oltre IVA..."
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("€ X.XXX,00 (XXXX/00)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Range($rng3.Start, $rng3.End).Font.Color = $wdColorRed
